$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before row 11 (shifts existing rows 11-18 down to 14-21)
$ws.Rows("11:13").Insert()

# Fill in the new enum entries: CD, SPEED, PIERCE continuing the AttributeType enum sequence
$ws.Range("G11").Value = "CD"
$ws.Range("I11").Value = 7

$ws.Range("G12").Value = "SPEED"
$ws.Range("I12").Value = 8

$ws.Range("G13").Value = "PIERCE"
$ws.Range("I13").Value = 9

# Match the saved cursor position recorded in the sheet view
$ws.Range("J13").Select() | Out-Null
